$p = $ppt.ActivePresentation

# Move slide 6 ("Besonderheit Security") to position 5, pushing
# slide 5 ("Besonderheit Systemanleitung") down to position 6.
$moveSlide = $p.Slides.Item(6)
$moveSlide.MoveTo(5)

# After the move, slide 5 is the former "Security" slide and
# slide 6 is the former "Systemanleitung" slide. Update their
# titles to add a " - " separator, matching the edited wording.
$securitySlide = $p.Slides.Item(5)
$securityTitle = $securitySlide.Shapes.Item(1).TextFrame.TextRange
$securityTitle.Text = "Besonderheit - "
[void]$securityTitle.InsertAfter("Security")

$systemanleitungSlide = $p.Slides.Item(6)
$systemanleitungTitle = $systemanleitungSlide.Shapes.Item(1).TextFrame.TextRange
$systemanleitungTitle.Text = "Besonderheit "
[void]$systemanleitungTitle.InsertAfter("- Systemanleitung")
